$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force the two new date-like labels to be written as plain text (matching the
# existing "dd-mm-yyyy" text entries already in column A) instead of being
# auto-converted to Excel date serials.
$ws.Range("A30:A31").NumberFormat = "@"
$ws.Range("A30").Value = "05-10-2021"
$ws.Range("A31").Value = "06-10-2021"
$ws.Range("A30:A31").Style = "Normal"

$ws.Range("B30").Value = 10000
$ws.Range("D30").Value = 0

$ws.Range("B31").Value = 10000
$ws.Range("D31").Value = 0
